# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.052.21'
$ws.Range('E2').Value = '  +2.32%  '

$ws.Range('D3').Value = '2.308.54'
$ws.Range('E3').Value = '  +2.01%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.61%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.68'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.38%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.73%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.513'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.09'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.94%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.87%  '

$ws.Range('E12').Value = '  +0.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.02'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.04%  '

$ws.Range('D14').Value = '2.665.83'
$ws.Range('E14').Value = '  +1.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.92'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.97%  '

$ws.Range('D16').Value = '2.306.35'
$ws.Range('E16').Value = '  +1.44%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.804'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.76%  '

$ws.Range('D18').Value = '43.025.01'
$ws.Range('E18').Value = '  +2.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.57'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.39%  '

$ws.Range('E20').Value = '  +2.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.21'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.75%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.46'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.18%  '

$ws.Range('E24').Value = '  +4.79%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.17'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.08'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.65'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.75%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.13'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.77'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.69%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.35'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.58%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.16'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.40%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.76'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.53%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0740'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.29%  '

$ws.Range('E37').Value = '  +0.31%  '

$ws.Range('E38').Value = '  +0.42%  '

$ws.Range('E39').Value = '  +1.17%  '

$ws.Range('E40').Value = '  +1.76%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.24'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.69%  '

$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.29'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.80%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.976.64'
$ws.Range('E43').Value = '  -0.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0289'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.09%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.29'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.91%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.01'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.80'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.98%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.98'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +18.52%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.12%  '

$ws.Range('E50').Value = '  +3.13%  '

$ws.Range('D51').Value = '2.533.40'
$ws.Range('E51').Value = '  +1.74%  '
